$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 16.4578
$ws.Range("B3").Value = 6.092200000000006
$ws.Range("B14").Value = 5.769800000000004
$ws.Range("B16").Value = 6.3266
$ws.Range("E18").Value = 18.18380000000001
$ws.Range("B21").Value = 8.785200000000005
$ws.Range("B23").Value = 8.510600000000005
$ws.Range("E24").Value = 16.2942
$ws.Range("B25").Value = 5.322099999999999
$ws.Range("E25").Value = 17.139
$ws.Range("B26").Value = 5.772800000000005
$ws.Range("E27").Value = 16.95789999999998
$ws.Range("B29").Value = 4.986700000000003
$ws.Range("E30").Value = 15.6255
$ws.Range("E31").Value = 16.40899999999999
$ws.Range("E39").Value = 16.0129
$ws.Range("B40").Value = 9.041399999999998
$ws.Range("E42").Value = 16.3999
$ws.Range("E48").Value = 17.3854
$ws.Range("E51").Value = 17.34270000000001
$ws.Range("E52").Value = 16.8601
$ws.Range("B53").Value = 4.845800000000001
$ws.Range("E55").Value = 16.46230000000001
$ws.Range("E56").Value = 16.3845
$ws.Range("B57").Value = 4.811599999999996
$ws.Range("E57").Value = 16.71760000000002
$ws.Range("B59").Value = 4.9025
$ws.Range("E60").Value = 16.2435
$ws.Range("B65").Value = 6.1698
$ws.Range("B69").Value = 5.203399999999999
$ws.Range("E73").Value = 17.58920000000001
$ws.Range("E74").Value = 16.83959999999998
$ws.Range("B79").Value = 8.986400000000005
$ws.Range("B83").Value = 5.639599999999999
$ws.Range("E89").Value = 17.40480000000001
$ws.Range("E90").Value = 16.90699999999999
$ws.Range("B91").Value = 5.2624
$ws.Range("E92").Value = 18.65540000000003
$ws.Range("B93").Value = 5.752000000000002
$ws.Range("B100").Value = 5.6144
